$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value (20000 -> 15000)
$ws.Range("C2").Value = 15000

# Fill in E3 (previously empty)
$ws.Range("E3").Value = "PERNR,ENDDA,BEGDA,MASSN,MASSG,STAT2"

# Fill rows 4-9 with new data (previously blank placeholder rows)
$data = @(
    @("employee", "S_PA0001", 20000, 3, "PERNR,ENDDA,BEGDA,BUKRS,WERKS,VDSK1,BTRTL,KOSTL,KOKRS,PERSG,PERSK,ORGEH,OTYPE,MSTBR"),
    @("employee", "S_PA0002", 20000, 4, "PERNR,ENDDA,BEGDA,INITS,NACHN,NACH2,VORNA,TITEL,MIDNM,SPRSL"),
    @("employee", "S_PA0006", 20000, 4, "SUBTY,ENDDA,BEGDA,ANSSA,STRAS,ORT01,ORT02,PSTLZ,LAND1,LOCAT,ADR03,ADR04,STATE,HSNMR,BLDNG,FLOOR,STRDS,COUNC,RCTVC,COM01,NUM01,COM02,NUM02,COM03,NUM03,COM04,NUM04,COM05,NUM05,COM06,NUM06"),
    @("employee", "S_PA0105", 20000, 5, "PERNR,SUBTY,ENDDA,BEGDA,USRTY"),
    @("employee", "S_INFOTYPE_TEXT", 20000, 6, "INFTY,SUBTY,ENDDA,BEGDA"),
    @("employee", "S_HRP1001", 20000, 7, "OTYPE,OBJID,BEGDA,ENDDA,SCLAS,SOBID,RSIGN,RELAT,PRIOX,PROZT")
)

$row = 4
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# Remove rows 11-17 (previously empty placeholder rows, no longer needed)
$ws.Range("A11:E17").EntireRow.Delete()
